# Insert two new rows before row 424 (rows 424 and 425 become new/blank;
# former rows 424-513 shift down to become rows 426-515).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("424:425").Insert() | Out-Null

# Fill in the new row 424 with fresh data.
$ws.Range("A424").Value2 = 7
$ws.Range("B424").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C424").Value2 = "Ñuble"
$ws.Range("D424").Value2 = 44785
$ws.Range("E424").Value2 = 16
$ws.Range("F424").Value2 = 100112020
$ws.Range("G424").Value2 = "Tomate"
$ws.Range("H424").Value2 = "Larga vida"
$ws.Range("I424").Value2 = "Primera"
$ws.Range("J424").Value2 = 300
$ws.Range("K424").Value2 = 7500
$ws.Range("L424").Value2 = 8000
$ws.Range("M424").Value2 = 7750
$ws.Range("N424").Value2 = "`$/bandeja 18 kilos"
$ws.Range("O424").Value2 = "Región de Arica y Parinacota"
$ws.Range("P424").Value2 = 431
$ws.Range("Q424").Value2 = 18
$ws.Range("R424").Value2 = "Hortaliza"

# Fill in the new row 425 with fresh data.
$ws.Range("A425").Value2 = 7
$ws.Range("B425").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C425").Value2 = "Ñuble"
$ws.Range("D425").Value2 = 44785
$ws.Range("E425").Value2 = 16
$ws.Range("F425").Value2 = 100112020
$ws.Range("G425").Value2 = "Tomate"
$ws.Range("H425").Value2 = "Larga vida"
$ws.Range("I425").Value2 = "Primera"
$ws.Range("J425").Value2 = 400
$ws.Range("K425").Value2 = 5000
$ws.Range("L425").Value2 = 5500
$ws.Range("M425").Value2 = 5250
$ws.Range("N425").Value2 = "`$/caja 10 kilos"
$ws.Range("O425").Value2 = "Región de Arica y Parinacota"
$ws.Range("P425").Value2 = 525
$ws.Range("Q425").Value2 = 10
$ws.Range("R425").Value2 = "Hortaliza"
